$d = $word.ActiveDocument

function Split-AtBoundary($paraIndex, $findText, $replaceText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) { throw "Could not find boundary '$findText' in paragraph $paraIndex" }
}

# --- Programa (Portuguese bullet list): insert a line break (w:br) before every bullet but the first ---
Split-AtBoundary 14 "- Reologia de fluidos,-" "- Reologia de fluidos,^l-"
Split-AtBoundary 14 "- Dimensionamento de tubulações,-" "- Dimensionamento de tubulações,^l-"
Split-AtBoundary 14 "- Acessórios e bombeamento para fluidos industriais-" "- Acessórios e bombeamento para fluidos industriais^l-"
Split-AtBoundary 14 "- Agitação e mistura-" "- Agitação e mistura^l-"
Split-AtBoundary 14 "- Caracterização de partículas e leito de partículas-" "- Caracterização de partículas e leito de partículas^l-"
Split-AtBoundary 14 "- Sedimentação-" "- Sedimentação^l-"
Split-AtBoundary 14 "- Filtração-" "- Filtração^l-"
Split-AtBoundary 14 "- Processos com membranas-" "- Processos com membranas^l-"
Split-AtBoundary 14 "- Operações unitárias de troca térmica: trocadores de calor e evaporadores-" "- Operações unitárias de troca térmica: trocadores de calor e evaporadores^l-"

# --- Programa (English bullet list): insert a line break (w:br) before every bullet but the first ---
Split-AtBoundary 15 "- Fluid rheology-" "- Fluid rheology^l-"
Split-AtBoundary 15 "- Sizing of pipes-" "- Sizing of pipes^l-"
Split-AtBoundary 15 "- Accessories and pumping for industrial fluids-" "- Accessories and pumping for industrial fluids^l-"
Split-AtBoundary 15 "- Stirring and mixing-" "- Stirring and mixing^l-"
Split-AtBoundary 15 "- Particle characterization and particle bed-" "- Particle characterization and particle bed^l-"
Split-AtBoundary 15 "- Sedimentation-" "- Sedimentation^l-"
Split-AtBoundary 15 "- Filtration-" "- Filtration^l-"
Split-AtBoundary 15 "- Processes with membranes-" "- Processes with membranes^l-"
Split-AtBoundary 15 "- Unit heat exchange operations: heat exchangers and evaporators-" "- Unit heat exchange operations: heat exchangers and evaporators^l-"

# --- Bibliografia (reference list): insert a line break (w:br) after every reference but the last ---
Split-AtBoundary 19 "FOUST, A.S., WENZEL, L. A., CLUMP, C.W., MAUS, L., ANDERSEN, L.B. Princípio das operações unitárias. Rio de Janeiro: Editora Guanabara Dois, 1982." "FOUST, A.S., WENZEL, L. A., CLUMP, C.W., MAUS, L., ANDERSEN, L.B. Princípio das operações unitárias. Rio de Janeiro: Editora Guanabara Dois, 1982.^l"
Split-AtBoundary 19 "GEANKOPLIS, C.J. Procesos de transporte y operaciones unitarias. Compañía Editorial Continental, S.A. de C.V. México, D.F., 1998." "GEANKOPLIS, C.J. Procesos de transporte y operaciones unitarias. Compañía Editorial Continental, S.A. de C.V. México, D.F., 1998.^l"
Split-AtBoundary 19 "PERRY, R.H. and CHILTON, C.H. Manual de Engenharia Química. 5a ed., Guanabara Dois, Rio de Janeiro, 1986." "PERRY, R.H. and CHILTON, C.H. Manual de Engenharia Química. 5a ed., Guanabara Dois, Rio de Janeiro, 1986.^l"
Split-AtBoundary 19 "REYNOLDS, T.D.; RICHARDS, P. Unit Operations and Processes in environmental Engineering. PWS Publishing, 1996." "REYNOLDS, T.D.; RICHARDS, P. Unit Operations and Processes in environmental Engineering. PWS Publishing, 1996.^l"

